# Commit: "Add functionality for commenting out lines in the xlsx forms to
# aid in form development"
#
# - survey sheet: turn the clause markers (if / end if / begin screen /
#   end screen) into commented-out markers (//if, // end if, //begin screen,
#   //end screen) and add "//" markers to the rows that were inside those
#   blocks, so the whole block is "commented out" for form development.
# - choices sheet: swap out thanksgiving/easter/halloween for
#   kwanzaa/christmas/hannukah/diwali (still 4 holiday choices, shifted).
# - survey/queries sheet tab selection moves back to the survey sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "survey" sheet - comment out the if/begin screen blocks
# ---------------------------------------------------------------------
$survey = $wb.Worksheets.Item("survey")

$survey.Range("B21").Value = "//if"
$survey.Range("B22").Value = "//"
$survey.Range("B23").Value = "// end if"

$survey.Range("B27").Value = "//begin screen"
$survey.Range("B28").Value = "//"
$survey.Range("B29").Value = "//"
$survey.Range("B30").Value = "//"
$survey.Range("B31").Value = "//"
$survey.Range("B32").Value = "//end screen"

# ---------------------------------------------------------------------
# "choices" sheet - replace thanksgiving/easter/halloween choices with
# kwanzaa/hannukah/diwali (and shift christmas down a row)
# ---------------------------------------------------------------------
$choices = $wb.Worksheets.Item("choices")

$choices.Range("B19").Value = "kwanzaa"
$choices.Range("D19").Value = "Kwanzaa"

$choices.Range("B20").Value = "christmas"
$choices.Range("D20").Value = "Christmas"

$choices.Range("B21").Value = "hannukah"
$choices.Range("D21").Value = "Hannukah"

$choices.Range("B22").Value = "diwali"
$choices.Range("D22").Value = "Diwali"

# ---------------------------------------------------------------------
# Selection / active-sheet bookkeeping to match the saved workbook state
# ---------------------------------------------------------------------
$queries = $wb.Worksheets.Item("queries")
$queries.Range("C7").Select() | Out-Null

$choices.Range("E20").Select() | Out-Null

$survey.Activate()
$survey.Range("B33").Select() | Out-Null
